$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.718.96"
$ws.Range("E2").Value = "  +0.53%  "
$ws.Range("D3").Value = "3.843.39"
$ws.Range("E3").Value = "  -1.59%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "522.96"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  +7.31%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.98"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  -2.98%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.605"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = "  -2.89%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.710"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = "  -4.55%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.169"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = "  -6.77%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000329"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = "  -7.41%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "41.54"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = "  -3.88%  "
$ws.Range("D13").Value = "4.479.37"
$ws.Range("E13").Value = "  -0.94%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.11"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "  -3.48%  "
$ws.Range("D15").Value = "3.869.65"
$ws.Range("E15").Value = "  +0.58%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.83"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = "  -2.87%  "
$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.134"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = "  -1.52%  "
$ws.Range("B18").Value = "Polygon"
$ws.Range("C18").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.20"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = "  +4.75%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "20.29"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = "  +1.30%  "
$ws.Range("D20").Value = "68.773.97"
$ws.Range("E20").Value = "  +0.47%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "419.61"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = "  -2.83%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.37"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = "  -5.29%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.04"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = "  -4.65%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "86.81"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  -3.13%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.95"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  +5.78%  "
$ws.Range("E26").Value = "  -7.41%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.47"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "  -4.55%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "35.88"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = "  -4.24%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "688.40"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = "  -3.52%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "13.06"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = "  -2.54%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.125"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = "  -4.70%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.83"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  -3.36%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "67.41"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = "  +9.47%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.431"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = "  +7.10%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.87"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = "  -3.53%  "
$ws.Range("D36").Value = "0.0₃0845"
$ws.Range("E36").Value = "  -6.05%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "39.56"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = "  -3.12%  "
$ws.Range("E38").Value = "  +0.07%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.146"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = "  -1.03%  "
$ws.Range("E40").Value = "  +0.15%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.22"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = "  +4.08%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0476"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = "  -4.08%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.15"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = "  +3.09%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.76"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = "  -6.76%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.38"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  +0.83%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.139"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = "  -2.77%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.95"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  +4.82%  "
$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").Value = "2.748.60"
$ws.Range("E48").Value = "  +14.09%  "
$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "144.03"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = "  +0.68%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.000267"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = "  +8.67%  "
$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D51").Value = "0.0₆0338"
$ws.Range("E51").Value = "  -10.84%  "
